$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference date-formatted cell style (column D uses a custom date/time numeric format)
$dateFormat = $ws.Cells.Item(2, 4).NumberFormat

# Row 201
$ws.Cells.Item(201, 1).Value = 11
$ws.Cells.Item(201, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(201, 3).Value = "Bíobío"
$ws.Cells.Item(201, 4).Value = 44609
$ws.Cells.Item(201, 5).Value = 8
$ws.Cells.Item(201, 6).Value = "Fruta"
$ws.Cells.Item(201, 7).Value = 100101
$ws.Cells.Item(201, 8).Value = "Berries"
$ws.Cells.Item(201, 9).Value = 100112025
$ws.Cells.Item(201, 10).Value = "Frutilla"
$ws.Cells.Item(201, 11).Value = "Sin especificar"
$ws.Cells.Item(201, 12).Value = "Especial"
$ws.Cells.Item(201, 13).Value = 100
$ws.Cells.Item(201, 14).Value = 7000
$ws.Cells.Item(201, 15).Value = 7000
$ws.Cells.Item(201, 16).Value = 7000
$ws.Cells.Item(201, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(201, 18).Value = "Quillón"
$ws.Cells.Item(201, 19).Value = 1000
$ws.Cells.Item(201, 20).Value = 7

# Row 202
$ws.Cells.Item(202, 1).Value = 11
$ws.Cells.Item(202, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(202, 3).Value = "Bíobío"
$ws.Cells.Item(202, 4).Value = 44609
$ws.Cells.Item(202, 5).Value = 8
$ws.Cells.Item(202, 6).Value = "Fruta"
$ws.Cells.Item(202, 7).Value = 100101
$ws.Cells.Item(202, 8).Value = "Berries"
$ws.Cells.Item(202, 9).Value = 100112025
$ws.Cells.Item(202, 10).Value = "Frutilla"
$ws.Cells.Item(202, 11).Value = "Sin especificar"
$ws.Cells.Item(202, 12).Value = "Primera"
$ws.Cells.Item(202, 13).Value = 100
$ws.Cells.Item(202, 14).Value = 6000
$ws.Cells.Item(202, 15).Value = 6000
$ws.Cells.Item(202, 16).Value = 6000
$ws.Cells.Item(202, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(202, 18).Value = "Quillón"
$ws.Cells.Item(202, 19).Value = 857
$ws.Cells.Item(202, 20).Value = 7

# Row 203
$ws.Cells.Item(203, 1).Value = 11
$ws.Cells.Item(203, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(203, 3).Value = "Bíobío"
$ws.Cells.Item(203, 4).Value = 44609
$ws.Cells.Item(203, 5).Value = 8
$ws.Cells.Item(203, 6).Value = "Fruta"
$ws.Cells.Item(203, 7).Value = 100101
$ws.Cells.Item(203, 8).Value = "Berries"
$ws.Cells.Item(203, 9).Value = 100112025
$ws.Cells.Item(203, 10).Value = "Frutilla"
$ws.Cells.Item(203, 11).Value = "Sin especificar"
$ws.Cells.Item(203, 12).Value = "Segunda"
$ws.Cells.Item(203, 13).Value = 100
$ws.Cells.Item(203, 14).Value = 5000
$ws.Cells.Item(203, 15).Value = 5000
$ws.Cells.Item(203, 16).Value = 5000
$ws.Cells.Item(203, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(203, 18).Value = "Quillón"
$ws.Cells.Item(203, 19).Value = 714
$ws.Cells.Item(203, 20).Value = 7

# Row 204
$ws.Cells.Item(204, 1).Value = 11
$ws.Cells.Item(204, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(204, 3).Value = "Bíobío"
$ws.Cells.Item(204, 4).Value = 44225
$ws.Cells.Item(204, 5).Value = 8
$ws.Cells.Item(204, 6).Value = "Fruta"
$ws.Cells.Item(204, 7).Value = 100101
$ws.Cells.Item(204, 8).Value = "Berries"
$ws.Cells.Item(204, 9).Value = 100112025
$ws.Cells.Item(204, 10).Value = "Frutilla"
$ws.Cells.Item(204, 11).Value = "Sin especificar"
$ws.Cells.Item(204, 12).Value = "Especial"
$ws.Cells.Item(204, 13).Value = 100
$ws.Cells.Item(204, 14).Value = 10000
$ws.Cells.Item(204, 15).Value = 10000
$ws.Cells.Item(204, 16).Value = 10000
$ws.Cells.Item(204, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(204, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(204, 19).Value = 1429
$ws.Cells.Item(204, 20).Value = 7

# Row 205
$ws.Cells.Item(205, 1).Value = 11
$ws.Cells.Item(205, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(205, 3).Value = "Bíobío"
$ws.Cells.Item(205, 4).Value = 44225
$ws.Cells.Item(205, 5).Value = 8
$ws.Cells.Item(205, 6).Value = "Fruta"
$ws.Cells.Item(205, 7).Value = 100101
$ws.Cells.Item(205, 8).Value = "Berries"
$ws.Cells.Item(205, 9).Value = 100112025
$ws.Cells.Item(205, 10).Value = "Frutilla"
$ws.Cells.Item(205, 11).Value = "Sin especificar"
$ws.Cells.Item(205, 12).Value = "Primera"
$ws.Cells.Item(205, 13).Value = 100
$ws.Cells.Item(205, 14).Value = 8000
$ws.Cells.Item(205, 15).Value = 8000
$ws.Cells.Item(205, 16).Value = 8000
$ws.Cells.Item(205, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(205, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(205, 19).Value = 1143
$ws.Cells.Item(205, 20).Value = 7

# Row 206
$ws.Cells.Item(206, 1).Value = 11
$ws.Cells.Item(206, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(206, 3).Value = "Bíobío"
$ws.Cells.Item(206, 4).Value = 44225
$ws.Cells.Item(206, 5).Value = 8
$ws.Cells.Item(206, 6).Value = "Fruta"
$ws.Cells.Item(206, 7).Value = 100101
$ws.Cells.Item(206, 8).Value = "Berries"
$ws.Cells.Item(206, 9).Value = 100112025
$ws.Cells.Item(206, 10).Value = "Frutilla"
$ws.Cells.Item(206, 11).Value = "Sin especificar"
$ws.Cells.Item(206, 12).Value = "Segunda"
$ws.Cells.Item(206, 13).Value = 50
$ws.Cells.Item(206, 14).Value = 6000
$ws.Cells.Item(206, 15).Value = 6000
$ws.Cells.Item(206, 16).Value = 6000
$ws.Cells.Item(206, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(206, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(206, 19).Value = 857
$ws.Cells.Item(206, 20).Value = 7

# Row 207
$ws.Cells.Item(207, 1).Value = 11
$ws.Cells.Item(207, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(207, 3).Value = "Bíobío"
$ws.Cells.Item(207, 4).Value = 44281
$ws.Cells.Item(207, 5).Value = 8
$ws.Cells.Item(207, 6).Value = "Fruta"
$ws.Cells.Item(207, 7).Value = 100101
$ws.Cells.Item(207, 8).Value = "Berries"
$ws.Cells.Item(207, 9).Value = 100112025
$ws.Cells.Item(207, 10).Value = "Frutilla"
$ws.Cells.Item(207, 11).Value = "Sin especificar"
$ws.Cells.Item(207, 12).Value = "Especial"
$ws.Cells.Item(207, 13).Value = 100
$ws.Cells.Item(207, 14).Value = 11000
$ws.Cells.Item(207, 15).Value = 11000
$ws.Cells.Item(207, 16).Value = 11000
$ws.Cells.Item(207, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(207, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(207, 19).Value = 1571
$ws.Cells.Item(207, 20).Value = 7

# Row 208
$ws.Cells.Item(208, 1).Value = 11
$ws.Cells.Item(208, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(208, 3).Value = "Bíobío"
$ws.Cells.Item(208, 4).Value = 44281
$ws.Cells.Item(208, 5).Value = 8
$ws.Cells.Item(208, 6).Value = "Fruta"
$ws.Cells.Item(208, 7).Value = 100101
$ws.Cells.Item(208, 8).Value = "Berries"
$ws.Cells.Item(208, 9).Value = 100112025
$ws.Cells.Item(208, 10).Value = "Frutilla"
$ws.Cells.Item(208, 11).Value = "Sin especificar"
$ws.Cells.Item(208, 12).Value = "Primera"
$ws.Cells.Item(208, 13).Value = 100
$ws.Cells.Item(208, 14).Value = 9000
$ws.Cells.Item(208, 15).Value = 9000
$ws.Cells.Item(208, 16).Value = 9000
$ws.Cells.Item(208, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(208, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(208, 19).Value = 1286
$ws.Cells.Item(208, 20).Value = 7

# Row 209
$ws.Cells.Item(209, 1).Value = 11
$ws.Cells.Item(209, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(209, 3).Value = "Bíobío"
$ws.Cells.Item(209, 4).Value = 44281
$ws.Cells.Item(209, 5).Value = 8
$ws.Cells.Item(209, 6).Value = "Fruta"
$ws.Cells.Item(209, 7).Value = 100101
$ws.Cells.Item(209, 8).Value = "Berries"
$ws.Cells.Item(209, 9).Value = 100112025
$ws.Cells.Item(209, 10).Value = "Frutilla"
$ws.Cells.Item(209, 11).Value = "Sin especificar"
$ws.Cells.Item(209, 12).Value = "Segunda"
$ws.Cells.Item(209, 13).Value = 50
$ws.Cells.Item(209, 14).Value = 7000
$ws.Cells.Item(209, 15).Value = 7000
$ws.Cells.Item(209, 16).Value = 7000
$ws.Cells.Item(209, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(209, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(209, 19).Value = 1000
$ws.Cells.Item(209, 20).Value = 7

# Row 210
$ws.Cells.Item(210, 1).Value = 11
$ws.Cells.Item(210, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(210, 3).Value = "Bíobío"
$ws.Cells.Item(210, 4).Value = 44252
$ws.Cells.Item(210, 5).Value = 8
$ws.Cells.Item(210, 6).Value = "Fruta"
$ws.Cells.Item(210, 7).Value = 100101
$ws.Cells.Item(210, 8).Value = "Berries"
$ws.Cells.Item(210, 9).Value = 100112025
$ws.Cells.Item(210, 10).Value = "Frutilla"
$ws.Cells.Item(210, 11).Value = "Sin especificar"
$ws.Cells.Item(210, 12).Value = "Especial"
$ws.Cells.Item(210, 13).Value = 100
$ws.Cells.Item(210, 14).Value = 10000
$ws.Cells.Item(210, 15).Value = 10000
$ws.Cells.Item(210, 16).Value = 10000
$ws.Cells.Item(210, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(210, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(210, 19).Value = 1429
$ws.Cells.Item(210, 20).Value = 7

# Row 211
$ws.Cells.Item(211, 1).Value = 11
$ws.Cells.Item(211, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(211, 3).Value = "Bíobío"
$ws.Cells.Item(211, 4).Value = 44252
$ws.Cells.Item(211, 5).Value = 8
$ws.Cells.Item(211, 6).Value = "Fruta"
$ws.Cells.Item(211, 7).Value = 100101
$ws.Cells.Item(211, 8).Value = "Berries"
$ws.Cells.Item(211, 9).Value = 100112025
$ws.Cells.Item(211, 10).Value = "Frutilla"
$ws.Cells.Item(211, 11).Value = "Sin especificar"
$ws.Cells.Item(211, 12).Value = "Primera"
$ws.Cells.Item(211, 13).Value = 100
$ws.Cells.Item(211, 14).Value = 8000
$ws.Cells.Item(211, 15).Value = 8000
$ws.Cells.Item(211, 16).Value = 8000
$ws.Cells.Item(211, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(211, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(211, 19).Value = 1143
$ws.Cells.Item(211, 20).Value = 7

# Row 212
$ws.Cells.Item(212, 1).Value = 11
$ws.Cells.Item(212, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(212, 3).Value = "Bíobío"
$ws.Cells.Item(212, 4).Value = 44252
$ws.Cells.Item(212, 5).Value = 8
$ws.Cells.Item(212, 6).Value = "Fruta"
$ws.Cells.Item(212, 7).Value = 100101
$ws.Cells.Item(212, 8).Value = "Berries"
$ws.Cells.Item(212, 9).Value = 100112025
$ws.Cells.Item(212, 10).Value = "Frutilla"
$ws.Cells.Item(212, 11).Value = "Sin especificar"
$ws.Cells.Item(212, 12).Value = "Segunda"
$ws.Cells.Item(212, 13).Value = 50
$ws.Cells.Item(212, 14).Value = 7000
$ws.Cells.Item(212, 15).Value = 7000
$ws.Cells.Item(212, 16).Value = 7000
$ws.Cells.Item(212, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(212, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(212, 19).Value = 1000
$ws.Cells.Item(212, 20).Value = 7

# Row 213
$ws.Cells.Item(213, 1).Value = 11
$ws.Cells.Item(213, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(213, 3).Value = "Bíobío"
$ws.Cells.Item(213, 4).Value = 44271
$ws.Cells.Item(213, 5).Value = 8
$ws.Cells.Item(213, 6).Value = "Fruta"
$ws.Cells.Item(213, 7).Value = 100101
$ws.Cells.Item(213, 8).Value = "Berries"
$ws.Cells.Item(213, 9).Value = 100112025
$ws.Cells.Item(213, 10).Value = "Frutilla"
$ws.Cells.Item(213, 11).Value = "Sin especificar"
$ws.Cells.Item(213, 12).Value = "Especial"
$ws.Cells.Item(213, 13).Value = 200
$ws.Cells.Item(213, 14).Value = 10000
$ws.Cells.Item(213, 15).Value = 10000
$ws.Cells.Item(213, 16).Value = 10000
$ws.Cells.Item(213, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(213, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(213, 19).Value = 1429
$ws.Cells.Item(213, 20).Value = 7

# Row 214
$ws.Cells.Item(214, 1).Value = 11
$ws.Cells.Item(214, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(214, 3).Value = "Bíobío"
$ws.Cells.Item(214, 4).Value = 44271
$ws.Cells.Item(214, 5).Value = 8
$ws.Cells.Item(214, 6).Value = "Fruta"
$ws.Cells.Item(214, 7).Value = 100101
$ws.Cells.Item(214, 8).Value = "Berries"
$ws.Cells.Item(214, 9).Value = 100112025
$ws.Cells.Item(214, 10).Value = "Frutilla"
$ws.Cells.Item(214, 11).Value = "Sin especificar"
$ws.Cells.Item(214, 12).Value = "Primera"
$ws.Cells.Item(214, 13).Value = 200
$ws.Cells.Item(214, 14).Value = 8000
$ws.Cells.Item(214, 15).Value = 8000
$ws.Cells.Item(214, 16).Value = 8000
$ws.Cells.Item(214, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(214, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(214, 19).Value = 1143
$ws.Cells.Item(214, 20).Value = 7

# Row 215
$ws.Cells.Item(215, 1).Value = 11
$ws.Cells.Item(215, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(215, 3).Value = "Bíobío"
$ws.Cells.Item(215, 4).Value = 44271
$ws.Cells.Item(215, 5).Value = 8
$ws.Cells.Item(215, 6).Value = "Fruta"
$ws.Cells.Item(215, 7).Value = 100101
$ws.Cells.Item(215, 8).Value = "Berries"
$ws.Cells.Item(215, 9).Value = 100112025
$ws.Cells.Item(215, 10).Value = "Frutilla"
$ws.Cells.Item(215, 11).Value = "Sin especificar"
$ws.Cells.Item(215, 12).Value = "Segunda"
$ws.Cells.Item(215, 13).Value = 100
$ws.Cells.Item(215, 14).Value = 6000
$ws.Cells.Item(215, 15).Value = 6000
$ws.Cells.Item(215, 16).Value = 6000
$ws.Cells.Item(215, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(215, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(215, 19).Value = 857
$ws.Cells.Item(215, 20).Value = 7

# Row 216
$ws.Cells.Item(216, 1).Value = 11
$ws.Cells.Item(216, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(216, 3).Value = "Bíobío"
$ws.Cells.Item(216, 4).Value = 44554
$ws.Cells.Item(216, 5).Value = 8
$ws.Cells.Item(216, 6).Value = "Fruta"
$ws.Cells.Item(216, 7).Value = 100101
$ws.Cells.Item(216, 8).Value = "Berries"
$ws.Cells.Item(216, 9).Value = 100112025
$ws.Cells.Item(216, 10).Value = "Frutilla"
$ws.Cells.Item(216, 11).Value = "Sin especificar"
$ws.Cells.Item(216, 12).Value = "Especial"
$ws.Cells.Item(216, 13).Value = 200
$ws.Cells.Item(216, 14).Value = 9000
$ws.Cells.Item(216, 15).Value = 9000
$ws.Cells.Item(216, 16).Value = 9000
$ws.Cells.Item(216, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(216, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(216, 19).Value = 1286
$ws.Cells.Item(216, 20).Value = 7

# Row 217
$ws.Cells.Item(217, 1).Value = 11
$ws.Cells.Item(217, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(217, 3).Value = "Bíobío"
$ws.Cells.Item(217, 4).Value = 44554
$ws.Cells.Item(217, 5).Value = 8
$ws.Cells.Item(217, 6).Value = "Fruta"
$ws.Cells.Item(217, 7).Value = 100101
$ws.Cells.Item(217, 8).Value = "Berries"
$ws.Cells.Item(217, 9).Value = 100112025
$ws.Cells.Item(217, 10).Value = "Frutilla"
$ws.Cells.Item(217, 11).Value = "Sin especificar"
$ws.Cells.Item(217, 12).Value = "Primera"
$ws.Cells.Item(217, 13).Value = 100
$ws.Cells.Item(217, 14).Value = 7000
$ws.Cells.Item(217, 15).Value = 7000
$ws.Cells.Item(217, 16).Value = 7000
$ws.Cells.Item(217, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(217, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(217, 19).Value = 1000
$ws.Cells.Item(217, 20).Value = 7

# Row 218
$ws.Cells.Item(218, 1).Value = 11
$ws.Cells.Item(218, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(218, 3).Value = "Bíobío"
$ws.Cells.Item(218, 4).Value = 44554
$ws.Cells.Item(218, 5).Value = 8
$ws.Cells.Item(218, 6).Value = "Fruta"
$ws.Cells.Item(218, 7).Value = 100101
$ws.Cells.Item(218, 8).Value = "Berries"
$ws.Cells.Item(218, 9).Value = 100112025
$ws.Cells.Item(218, 10).Value = "Frutilla"
$ws.Cells.Item(218, 11).Value = "Sin especificar"
$ws.Cells.Item(218, 12).Value = "Segunda"
$ws.Cells.Item(218, 13).Value = 100
$ws.Cells.Item(218, 14).Value = 6000
$ws.Cells.Item(218, 15).Value = 6000
$ws.Cells.Item(218, 16).Value = 6000
$ws.Cells.Item(218, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(218, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(218, 19).Value = 857
$ws.Cells.Item(218, 20).Value = 7

# Row 219
$ws.Cells.Item(219, 1).Value = 11
$ws.Cells.Item(219, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(219, 3).Value = "Bíobío"
$ws.Cells.Item(219, 4).Value = 44162
$ws.Cells.Item(219, 5).Value = 8
$ws.Cells.Item(219, 6).Value = "Fruta"
$ws.Cells.Item(219, 7).Value = 100101
$ws.Cells.Item(219, 8).Value = "Berries"
$ws.Cells.Item(219, 9).Value = 100112025
$ws.Cells.Item(219, 10).Value = "Frutilla"
$ws.Cells.Item(219, 11).Value = "Sin especificar"
$ws.Cells.Item(219, 12).Value = "Especial"
$ws.Cells.Item(219, 13).Value = 150
$ws.Cells.Item(219, 14).Value = 11000
$ws.Cells.Item(219, 15).Value = 11000
$ws.Cells.Item(219, 16).Value = 11000
$ws.Cells.Item(219, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(219, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(219, 19).Value = 1571
$ws.Cells.Item(219, 20).Value = 7

# Row 220
$ws.Cells.Item(220, 1).Value = 11
$ws.Cells.Item(220, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(220, 3).Value = "Bíobío"
$ws.Cells.Item(220, 4).Value = 44162
$ws.Cells.Item(220, 5).Value = 8
$ws.Cells.Item(220, 6).Value = "Fruta"
$ws.Cells.Item(220, 7).Value = 100101
$ws.Cells.Item(220, 8).Value = "Berries"
$ws.Cells.Item(220, 9).Value = 100112025
$ws.Cells.Item(220, 10).Value = "Frutilla"
$ws.Cells.Item(220, 11).Value = "Sin especificar"
$ws.Cells.Item(220, 12).Value = "Primera"
$ws.Cells.Item(220, 13).Value = 150
$ws.Cells.Item(220, 14).Value = 9000
$ws.Cells.Item(220, 15).Value = 9000
$ws.Cells.Item(220, 16).Value = 9000
$ws.Cells.Item(220, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(220, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(220, 19).Value = 1286
$ws.Cells.Item(220, 20).Value = 7

# Row 221
$ws.Cells.Item(221, 1).Value = 11
$ws.Cells.Item(221, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(221, 3).Value = "Bíobío"
$ws.Cells.Item(221, 4).Value = 44162
$ws.Cells.Item(221, 5).Value = 8
$ws.Cells.Item(221, 6).Value = "Fruta"
$ws.Cells.Item(221, 7).Value = 100101
$ws.Cells.Item(221, 8).Value = "Berries"
$ws.Cells.Item(221, 9).Value = 100112025
$ws.Cells.Item(221, 10).Value = "Frutilla"
$ws.Cells.Item(221, 11).Value = "Sin especificar"
$ws.Cells.Item(221, 12).Value = "Segunda"
$ws.Cells.Item(221, 13).Value = 100
$ws.Cells.Item(221, 14).Value = 7000
$ws.Cells.Item(221, 15).Value = 7000
$ws.Cells.Item(221, 16).Value = 7000
$ws.Cells.Item(221, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(221, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(221, 19).Value = 1000
$ws.Cells.Item(221, 20).Value = 7

# Row 222
$ws.Cells.Item(222, 1).Value = 11
$ws.Cells.Item(222, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(222, 3).Value = "Bíobío"
$ws.Cells.Item(222, 4).Value = 44516
$ws.Cells.Item(222, 5).Value = 8
$ws.Cells.Item(222, 6).Value = "Fruta"
$ws.Cells.Item(222, 7).Value = 100101
$ws.Cells.Item(222, 8).Value = "Berries"
$ws.Cells.Item(222, 9).Value = 100112025
$ws.Cells.Item(222, 10).Value = "Frutilla"
$ws.Cells.Item(222, 11).Value = "Sin especificar"
$ws.Cells.Item(222, 12).Value = "Especial"
$ws.Cells.Item(222, 13).Value = 200
$ws.Cells.Item(222, 14).Value = 8000
$ws.Cells.Item(222, 15).Value = 8000
$ws.Cells.Item(222, 16).Value = 8000
$ws.Cells.Item(222, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(222, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(222, 19).Value = 1143
$ws.Cells.Item(222, 20).Value = 7

# Row 223
$ws.Cells.Item(223, 1).Value = 11
$ws.Cells.Item(223, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(223, 3).Value = "Bíobío"
$ws.Cells.Item(223, 4).Value = 44516
$ws.Cells.Item(223, 5).Value = 8
$ws.Cells.Item(223, 6).Value = "Fruta"
$ws.Cells.Item(223, 7).Value = 100101
$ws.Cells.Item(223, 8).Value = "Berries"
$ws.Cells.Item(223, 9).Value = 100112025
$ws.Cells.Item(223, 10).Value = "Frutilla"
$ws.Cells.Item(223, 11).Value = "Sin especificar"
$ws.Cells.Item(223, 12).Value = "Primera"
$ws.Cells.Item(223, 13).Value = 1500
$ws.Cells.Item(223, 14).Value = 6500
$ws.Cells.Item(223, 15).Value = 6500
$ws.Cells.Item(223, 16).Value = 6500
$ws.Cells.Item(223, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(223, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(223, 19).Value = 929
$ws.Cells.Item(223, 20).Value = 7

# Row 224
$ws.Cells.Item(224, 1).Value = 11
$ws.Cells.Item(224, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(224, 3).Value = "Bíobío"
$ws.Cells.Item(224, 4).Value = 44516
$ws.Cells.Item(224, 5).Value = 8
$ws.Cells.Item(224, 6).Value = "Fruta"
$ws.Cells.Item(224, 7).Value = 100101
$ws.Cells.Item(224, 8).Value = "Berries"
$ws.Cells.Item(224, 9).Value = 100112025
$ws.Cells.Item(224, 10).Value = "Frutilla"
$ws.Cells.Item(224, 11).Value = "Sin especificar"
$ws.Cells.Item(224, 12).Value = "Segunda"
$ws.Cells.Item(224, 13).Value = 100
$ws.Cells.Item(224, 14).Value = 5000
$ws.Cells.Item(224, 15).Value = 5000
$ws.Cells.Item(224, 16).Value = 5000
$ws.Cells.Item(224, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(224, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(224, 19).Value = 714
$ws.Cells.Item(224, 20).Value = 7

# Row 225
$ws.Cells.Item(225, 1).Value = 11
$ws.Cells.Item(225, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(225, 3).Value = "Bíobío"
$ws.Cells.Item(225, 4).Value = 44511
$ws.Cells.Item(225, 5).Value = 8
$ws.Cells.Item(225, 6).Value = "Fruta"
$ws.Cells.Item(225, 7).Value = 100101
$ws.Cells.Item(225, 8).Value = "Berries"
$ws.Cells.Item(225, 9).Value = 100112025
$ws.Cells.Item(225, 10).Value = "Frutilla"
$ws.Cells.Item(225, 11).Value = "Sin especificar"
$ws.Cells.Item(225, 12).Value = "Primera"
$ws.Cells.Item(225, 13).Value = 600
$ws.Cells.Item(225, 14).Value = 7500
$ws.Cells.Item(225, 15).Value = 8000
$ws.Cells.Item(225, 16).Value = 7792
$ws.Cells.Item(225, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(225, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(225, 19).Value = 1113
$ws.Cells.Item(225, 20).Value = 7

# Row 226
$ws.Cells.Item(226, 1).Value = 11
$ws.Cells.Item(226, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(226, 3).Value = "Bíobío"
$ws.Cells.Item(226, 4).Value = 44511
$ws.Cells.Item(226, 5).Value = 8
$ws.Cells.Item(226, 6).Value = "Fruta"
$ws.Cells.Item(226, 7).Value = 100101
$ws.Cells.Item(226, 8).Value = "Berries"
$ws.Cells.Item(226, 9).Value = 100112025
$ws.Cells.Item(226, 10).Value = "Frutilla"
$ws.Cells.Item(226, 11).Value = "Sin especificar"
$ws.Cells.Item(226, 12).Value = "Segunda"
$ws.Cells.Item(226, 13).Value = 450
$ws.Cells.Item(226, 14).Value = 6000
$ws.Cells.Item(226, 15).Value = 6500
$ws.Cells.Item(226, 16).Value = 6222
$ws.Cells.Item(226, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(226, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(226, 19).Value = 889
$ws.Cells.Item(226, 20).Value = 7

# Row 227
$ws.Cells.Item(227, 1).Value = 11
$ws.Cells.Item(227, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(227, 3).Value = "Bíobío"
$ws.Cells.Item(227, 4).Value = 44565
$ws.Cells.Item(227, 5).Value = 8
$ws.Cells.Item(227, 6).Value = "Fruta"
$ws.Cells.Item(227, 7).Value = 100101
$ws.Cells.Item(227, 8).Value = "Berries"
$ws.Cells.Item(227, 9).Value = 100112025
$ws.Cells.Item(227, 10).Value = "Frutilla"
$ws.Cells.Item(227, 11).Value = "Sin especificar"
$ws.Cells.Item(227, 12).Value = "Primera"
$ws.Cells.Item(227, 13).Value = 250
$ws.Cells.Item(227, 14).Value = 5500
$ws.Cells.Item(227, 15).Value = 6000
$ws.Cells.Item(227, 16).Value = 5800
$ws.Cells.Item(227, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(227, 18).Value = "Región del Maule"
$ws.Cells.Item(227, 19).Value = 829
$ws.Cells.Item(227, 20).Value = 7

# Row 228
$ws.Cells.Item(228, 1).Value = 11
$ws.Cells.Item(228, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(228, 3).Value = "Bíobío"
$ws.Cells.Item(228, 4).Value = 44573
$ws.Cells.Item(228, 5).Value = 8
$ws.Cells.Item(228, 6).Value = "Fruta"
$ws.Cells.Item(228, 7).Value = 100101
$ws.Cells.Item(228, 8).Value = "Berries"
$ws.Cells.Item(228, 9).Value = 100112025
$ws.Cells.Item(228, 10).Value = "Frutilla"
$ws.Cells.Item(228, 11).Value = "Sin especificar"
$ws.Cells.Item(228, 12).Value = "Especial"
$ws.Cells.Item(228, 13).Value = 200
$ws.Cells.Item(228, 14).Value = 8000
$ws.Cells.Item(228, 15).Value = 8000
$ws.Cells.Item(228, 16).Value = 8000
$ws.Cells.Item(228, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(228, 18).Value = "Quillón"
$ws.Cells.Item(228, 19).Value = 1143
$ws.Cells.Item(228, 20).Value = 7

# Row 229
$ws.Cells.Item(229, 1).Value = 11
$ws.Cells.Item(229, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(229, 3).Value = "Bíobío"
$ws.Cells.Item(229, 4).Value = 44573
$ws.Cells.Item(229, 5).Value = 8
$ws.Cells.Item(229, 6).Value = "Fruta"
$ws.Cells.Item(229, 7).Value = 100101
$ws.Cells.Item(229, 8).Value = "Berries"
$ws.Cells.Item(229, 9).Value = 100112025
$ws.Cells.Item(229, 10).Value = "Frutilla"
$ws.Cells.Item(229, 11).Value = "Sin especificar"
$ws.Cells.Item(229, 12).Value = "Primera"
$ws.Cells.Item(229, 13).Value = 200
$ws.Cells.Item(229, 14).Value = 6000
$ws.Cells.Item(229, 15).Value = 6000
$ws.Cells.Item(229, 16).Value = 6000
$ws.Cells.Item(229, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(229, 18).Value = "Quillón"
$ws.Cells.Item(229, 19).Value = 857
$ws.Cells.Item(229, 20).Value = 7

# Row 230
$ws.Cells.Item(230, 1).Value = 11
$ws.Cells.Item(230, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(230, 3).Value = "Bíobío"
$ws.Cells.Item(230, 4).Value = 44573
$ws.Cells.Item(230, 5).Value = 8
$ws.Cells.Item(230, 6).Value = "Fruta"
$ws.Cells.Item(230, 7).Value = 100101
$ws.Cells.Item(230, 8).Value = "Berries"
$ws.Cells.Item(230, 9).Value = 100112025
$ws.Cells.Item(230, 10).Value = "Frutilla"
$ws.Cells.Item(230, 11).Value = "Sin especificar"
$ws.Cells.Item(230, 12).Value = "Segunda"
$ws.Cells.Item(230, 13).Value = 100
$ws.Cells.Item(230, 14).Value = 5000
$ws.Cells.Item(230, 15).Value = 5000
$ws.Cells.Item(230, 16).Value = 5000
$ws.Cells.Item(230, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(230, 18).Value = "Quillón"
$ws.Cells.Item(230, 19).Value = 714
$ws.Cells.Item(230, 20).Value = 7

# Row 231
$ws.Cells.Item(231, 1).Value = 11
$ws.Cells.Item(231, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(231, 3).Value = "Bíobío"
$ws.Cells.Item(231, 4).Value = 44567
$ws.Cells.Item(231, 5).Value = 8
$ws.Cells.Item(231, 6).Value = "Fruta"
$ws.Cells.Item(231, 7).Value = 100101
$ws.Cells.Item(231, 8).Value = "Berries"
$ws.Cells.Item(231, 9).Value = 100112025
$ws.Cells.Item(231, 10).Value = "Frutilla"
$ws.Cells.Item(231, 11).Value = "Sin especificar"
$ws.Cells.Item(231, 12).Value = "Primera"
$ws.Cells.Item(231, 13).Value = 250
$ws.Cells.Item(231, 14).Value = 7000
$ws.Cells.Item(231, 15).Value = 7500
$ws.Cells.Item(231, 16).Value = 7260
$ws.Cells.Item(231, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(231, 18).Value = "Región del Maule"
$ws.Cells.Item(231, 19).Value = 1037
$ws.Cells.Item(231, 20).Value = 7

# Row 232
$ws.Cells.Item(232, 1).Value = 11
$ws.Cells.Item(232, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(232, 3).Value = "Bíobío"
$ws.Cells.Item(232, 4).Value = 44567
$ws.Cells.Item(232, 5).Value = 8
$ws.Cells.Item(232, 6).Value = "Fruta"
$ws.Cells.Item(232, 7).Value = 100101
$ws.Cells.Item(232, 8).Value = "Berries"
$ws.Cells.Item(232, 9).Value = 100112025
$ws.Cells.Item(232, 10).Value = "Frutilla"
$ws.Cells.Item(232, 11).Value = "Sin especificar"
$ws.Cells.Item(232, 12).Value = "Segunda"
$ws.Cells.Item(232, 13).Value = 270
$ws.Cells.Item(232, 14).Value = 6000
$ws.Cells.Item(232, 15).Value = 6500
$ws.Cells.Item(232, 16).Value = 6259
$ws.Cells.Item(232, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(232, 18).Value = "Región del Maule"
$ws.Cells.Item(232, 19).Value = 894
$ws.Cells.Item(232, 20).Value = 7

# Row 233
$ws.Cells.Item(233, 1).Value = 11
$ws.Cells.Item(233, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(233, 3).Value = "Bíobío"
$ws.Cells.Item(233, 4).Value = 44581
$ws.Cells.Item(233, 5).Value = 8
$ws.Cells.Item(233, 6).Value = "Fruta"
$ws.Cells.Item(233, 7).Value = 100101
$ws.Cells.Item(233, 8).Value = "Berries"
$ws.Cells.Item(233, 9).Value = 100112025
$ws.Cells.Item(233, 10).Value = "Frutilla"
$ws.Cells.Item(233, 11).Value = "Sin especificar"
$ws.Cells.Item(233, 12).Value = "Especial"
$ws.Cells.Item(233, 13).Value = 280
$ws.Cells.Item(233, 14).Value = 7000
$ws.Cells.Item(233, 15).Value = 7500
$ws.Cells.Item(233, 16).Value = 7232
$ws.Cells.Item(233, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(233, 18).Value = "Región del Maule"
$ws.Cells.Item(233, 19).Value = 1033
$ws.Cells.Item(233, 20).Value = 7

# Row 234
$ws.Cells.Item(234, 1).Value = 11
$ws.Cells.Item(234, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(234, 3).Value = "Bíobío"
$ws.Cells.Item(234, 4).Value = 44581
$ws.Cells.Item(234, 5).Value = 8
$ws.Cells.Item(234, 6).Value = "Fruta"
$ws.Cells.Item(234, 7).Value = 100101
$ws.Cells.Item(234, 8).Value = "Berries"
$ws.Cells.Item(234, 9).Value = 100112025
$ws.Cells.Item(234, 10).Value = "Frutilla"
$ws.Cells.Item(234, 11).Value = "Sin especificar"
$ws.Cells.Item(234, 12).Value = "Primera"
$ws.Cells.Item(234, 13).Value = 250
$ws.Cells.Item(234, 14).Value = 6000
$ws.Cells.Item(234, 15).Value = 6500
$ws.Cells.Item(234, 16).Value = 6260
$ws.Cells.Item(234, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(234, 18).Value = "Región del Maule"
$ws.Cells.Item(234, 19).Value = 894
$ws.Cells.Item(234, 20).Value = 7

# Row 235
$ws.Cells.Item(235, 1).Value = 11
$ws.Cells.Item(235, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(235, 3).Value = "Bíobío"
$ws.Cells.Item(235, 4).Value = 44537
$ws.Cells.Item(235, 5).Value = 8
$ws.Cells.Item(235, 6).Value = "Fruta"
$ws.Cells.Item(235, 7).Value = 100101
$ws.Cells.Item(235, 8).Value = "Berries"
$ws.Cells.Item(235, 9).Value = 100112025
$ws.Cells.Item(235, 10).Value = "Frutilla"
$ws.Cells.Item(235, 11).Value = "Sin especificar"
$ws.Cells.Item(235, 12).Value = "Especial"
$ws.Cells.Item(235, 13).Value = 250
$ws.Cells.Item(235, 14).Value = 8000
$ws.Cells.Item(235, 15).Value = 9500
$ws.Cells.Item(235, 16).Value = 8600
$ws.Cells.Item(235, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(235, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(235, 19).Value = 1229
$ws.Cells.Item(235, 20).Value = 7

# Row 236
$ws.Cells.Item(236, 1).Value = 11
$ws.Cells.Item(236, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(236, 3).Value = "Bíobío"
$ws.Cells.Item(236, 4).Value = 44537
$ws.Cells.Item(236, 5).Value = 8
$ws.Cells.Item(236, 6).Value = "Fruta"
$ws.Cells.Item(236, 7).Value = 100101
$ws.Cells.Item(236, 8).Value = "Berries"
$ws.Cells.Item(236, 9).Value = 100112025
$ws.Cells.Item(236, 10).Value = "Frutilla"
$ws.Cells.Item(236, 11).Value = "Sin especificar"
$ws.Cells.Item(236, 12).Value = "Primera"
$ws.Cells.Item(236, 13).Value = 220
$ws.Cells.Item(236, 14).Value = 5500
$ws.Cells.Item(236, 15).Value = 6000
$ws.Cells.Item(236, 16).Value = 5773
$ws.Cells.Item(236, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(236, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(236, 19).Value = 825
$ws.Cells.Item(236, 20).Value = 7

# Row 237
$ws.Cells.Item(237, 1).Value = 11
$ws.Cells.Item(237, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(237, 3).Value = "Bíobío"
$ws.Cells.Item(237, 4).Value = 44553
$ws.Cells.Item(237, 5).Value = 8
$ws.Cells.Item(237, 6).Value = "Fruta"
$ws.Cells.Item(237, 7).Value = 100101
$ws.Cells.Item(237, 8).Value = "Berries"
$ws.Cells.Item(237, 9).Value = 100112025
$ws.Cells.Item(237, 10).Value = "Frutilla"
$ws.Cells.Item(237, 11).Value = "Sin especificar"
$ws.Cells.Item(237, 12).Value = "Primera"
$ws.Cells.Item(237, 13).Value = 200
$ws.Cells.Item(237, 14).Value = 7000
$ws.Cells.Item(237, 15).Value = 7500
$ws.Cells.Item(237, 16).Value = 7250
$ws.Cells.Item(237, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(237, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(237, 19).Value = 1036
$ws.Cells.Item(237, 20).Value = 7

# Row 238
$ws.Cells.Item(238, 1).Value = 11
$ws.Cells.Item(238, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(238, 3).Value = "Bíobío"
$ws.Cells.Item(238, 4).Value = 44490
$ws.Cells.Item(238, 5).Value = 8
$ws.Cells.Item(238, 6).Value = "Fruta"
$ws.Cells.Item(238, 7).Value = 100101
$ws.Cells.Item(238, 8).Value = "Berries"
$ws.Cells.Item(238, 9).Value = 100112025
$ws.Cells.Item(238, 10).Value = "Frutilla"
$ws.Cells.Item(238, 11).Value = "Sin especificar"
$ws.Cells.Item(238, 12).Value = "Primera"
$ws.Cells.Item(238, 13).Value = 450
$ws.Cells.Item(238, 14).Value = 8000
$ws.Cells.Item(238, 15).Value = 8500
$ws.Cells.Item(238, 16).Value = 8222
$ws.Cells.Item(238, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(238, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(238, 19).Value = 1175
$ws.Cells.Item(238, 20).Value = 7

# Row 239
$ws.Cells.Item(239, 1).Value = 11
$ws.Cells.Item(239, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(239, 3).Value = "Bíobío"
$ws.Cells.Item(239, 4).Value = 44490
$ws.Cells.Item(239, 5).Value = 8
$ws.Cells.Item(239, 6).Value = "Fruta"
$ws.Cells.Item(239, 7).Value = 100101
$ws.Cells.Item(239, 8).Value = "Berries"
$ws.Cells.Item(239, 9).Value = 100112025
$ws.Cells.Item(239, 10).Value = "Frutilla"
$ws.Cells.Item(239, 11).Value = "Sin especificar"
$ws.Cells.Item(239, 12).Value = "Segunda"
$ws.Cells.Item(239, 13).Value = 430
$ws.Cells.Item(239, 14).Value = 7000
$ws.Cells.Item(239, 15).Value = 7500
$ws.Cells.Item(239, 16).Value = 7233
$ws.Cells.Item(239, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(239, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(239, 19).Value = 1033
$ws.Cells.Item(239, 20).Value = 7

# Row 240
$ws.Cells.Item(240, 1).Value = 11
$ws.Cells.Item(240, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(240, 3).Value = "Bíobío"
$ws.Cells.Item(240, 4).Value = 44446
$ws.Cells.Item(240, 5).Value = 8
$ws.Cells.Item(240, 6).Value = "Fruta"
$ws.Cells.Item(240, 7).Value = 100101
$ws.Cells.Item(240, 8).Value = "Berries"
$ws.Cells.Item(240, 9).Value = 100112025
$ws.Cells.Item(240, 10).Value = "Frutilla"
$ws.Cells.Item(240, 11).Value = "Sin especificar"
$ws.Cells.Item(240, 12).Value = "Primera"
$ws.Cells.Item(240, 13).Value = 100
$ws.Cells.Item(240, 14).Value = 21000
$ws.Cells.Item(240, 15).Value = 22000
$ws.Cells.Item(240, 16).Value = 21500
$ws.Cells.Item(240, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(240, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(240, 19).Value = 3071
$ws.Cells.Item(240, 20).Value = 7

# Row 241
$ws.Cells.Item(241, 1).Value = 11
$ws.Cells.Item(241, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(241, 3).Value = "Bíobío"
$ws.Cells.Item(241, 4).Value = 44421
$ws.Cells.Item(241, 5).Value = 8
$ws.Cells.Item(241, 6).Value = "Fruta"
$ws.Cells.Item(241, 7).Value = 100101
$ws.Cells.Item(241, 8).Value = "Berries"
$ws.Cells.Item(241, 9).Value = 100112025
$ws.Cells.Item(241, 10).Value = "Frutilla"
$ws.Cells.Item(241, 11).Value = "Sin especificar"
$ws.Cells.Item(241, 12).Value = "Primera"
$ws.Cells.Item(241, 13).Value = 50
$ws.Cells.Item(241, 14).Value = 22000
$ws.Cells.Item(241, 15).Value = 22000
$ws.Cells.Item(241, 16).Value = 22000
$ws.Cells.Item(241, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(241, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(241, 19).Value = 3143
$ws.Cells.Item(241, 20).Value = 7

# Row 242
$ws.Cells.Item(242, 1).Value = 11
$ws.Cells.Item(242, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(242, 3).Value = "Bíobío"
$ws.Cells.Item(242, 4).Value = 44208
$ws.Cells.Item(242, 5).Value = 8
$ws.Cells.Item(242, 6).Value = "Fruta"
$ws.Cells.Item(242, 7).Value = 100101
$ws.Cells.Item(242, 8).Value = "Berries"
$ws.Cells.Item(242, 9).Value = 100112025
$ws.Cells.Item(242, 10).Value = "Frutilla"
$ws.Cells.Item(242, 11).Value = "Sin especificar"
$ws.Cells.Item(242, 12).Value = "Especial"
$ws.Cells.Item(242, 13).Value = 200
$ws.Cells.Item(242, 14).Value = 10000
$ws.Cells.Item(242, 15).Value = 10000
$ws.Cells.Item(242, 16).Value = 10000
$ws.Cells.Item(242, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(242, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(242, 19).Value = 1429
$ws.Cells.Item(242, 20).Value = 7

# Row 243
$ws.Cells.Item(243, 1).Value = 11
$ws.Cells.Item(243, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(243, 3).Value = "Bíobío"
$ws.Cells.Item(243, 4).Value = 44208
$ws.Cells.Item(243, 5).Value = 8
$ws.Cells.Item(243, 6).Value = "Fruta"
$ws.Cells.Item(243, 7).Value = 100101
$ws.Cells.Item(243, 8).Value = "Berries"
$ws.Cells.Item(243, 9).Value = 100112025
$ws.Cells.Item(243, 10).Value = "Frutilla"
$ws.Cells.Item(243, 11).Value = "Sin especificar"
$ws.Cells.Item(243, 12).Value = "Primera"
$ws.Cells.Item(243, 13).Value = 150
$ws.Cells.Item(243, 14).Value = 8000
$ws.Cells.Item(243, 15).Value = 8000
$ws.Cells.Item(243, 16).Value = 8000
$ws.Cells.Item(243, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(243, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(243, 19).Value = 1143
$ws.Cells.Item(243, 20).Value = 7

# Row 244
$ws.Cells.Item(244, 1).Value = 11
$ws.Cells.Item(244, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(244, 3).Value = "Bíobío"
$ws.Cells.Item(244, 4).Value = 44208
$ws.Cells.Item(244, 5).Value = 8
$ws.Cells.Item(244, 6).Value = "Fruta"
$ws.Cells.Item(244, 7).Value = 100101
$ws.Cells.Item(244, 8).Value = "Berries"
$ws.Cells.Item(244, 9).Value = 100112025
$ws.Cells.Item(244, 10).Value = "Frutilla"
$ws.Cells.Item(244, 11).Value = "Sin especificar"
$ws.Cells.Item(244, 12).Value = "Segunda"
$ws.Cells.Item(244, 13).Value = 50
$ws.Cells.Item(244, 14).Value = 7000
$ws.Cells.Item(244, 15).Value = 7000
$ws.Cells.Item(244, 16).Value = 7000
$ws.Cells.Item(244, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(244, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(244, 19).Value = 1000
$ws.Cells.Item(244, 20).Value = 7

# Row 245
$ws.Cells.Item(245, 1).Value = 11
$ws.Cells.Item(245, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(245, 3).Value = "Bíobío"
$ws.Cells.Item(245, 4).Value = 44355
$ws.Cells.Item(245, 5).Value = 8
$ws.Cells.Item(245, 6).Value = "Fruta"
$ws.Cells.Item(245, 7).Value = 100101
$ws.Cells.Item(245, 8).Value = "Berries"
$ws.Cells.Item(245, 9).Value = 100112025
$ws.Cells.Item(245, 10).Value = "Frutilla"
$ws.Cells.Item(245, 11).Value = "Sin especificar"
$ws.Cells.Item(245, 12).Value = "Especial"
$ws.Cells.Item(245, 13).Value = 50
$ws.Cells.Item(245, 14).Value = 14000
$ws.Cells.Item(245, 15).Value = 14000
$ws.Cells.Item(245, 16).Value = 14000
$ws.Cells.Item(245, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(245, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(245, 19).Value = 2000
$ws.Cells.Item(245, 20).Value = 7

# Row 246
$ws.Cells.Item(246, 1).Value = 11
$ws.Cells.Item(246, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(246, 3).Value = "Bíobío"
$ws.Cells.Item(246, 4).Value = 44355
$ws.Cells.Item(246, 5).Value = 8
$ws.Cells.Item(246, 6).Value = "Fruta"
$ws.Cells.Item(246, 7).Value = 100101
$ws.Cells.Item(246, 8).Value = "Berries"
$ws.Cells.Item(246, 9).Value = 100112025
$ws.Cells.Item(246, 10).Value = "Frutilla"
$ws.Cells.Item(246, 11).Value = "Sin especificar"
$ws.Cells.Item(246, 12).Value = "Primera"
$ws.Cells.Item(246, 13).Value = 50
$ws.Cells.Item(246, 14).Value = 11000
$ws.Cells.Item(246, 15).Value = 11000
$ws.Cells.Item(246, 16).Value = 11000
$ws.Cells.Item(246, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(246, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(246, 19).Value = 1571
$ws.Cells.Item(246, 20).Value = 7

# Row 247
$ws.Cells.Item(247, 1).Value = 11
$ws.Cells.Item(247, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(247, 3).Value = "Bíobío"
$ws.Cells.Item(247, 4).Value = 44530
$ws.Cells.Item(247, 5).Value = 8
$ws.Cells.Item(247, 6).Value = "Fruta"
$ws.Cells.Item(247, 7).Value = 100101
$ws.Cells.Item(247, 8).Value = "Berries"
$ws.Cells.Item(247, 9).Value = 100112025
$ws.Cells.Item(247, 10).Value = "Frutilla"
$ws.Cells.Item(247, 11).Value = "Sin especificar"
$ws.Cells.Item(247, 12).Value = "Especial"
$ws.Cells.Item(247, 13).Value = 200
$ws.Cells.Item(247, 14).Value = 9000
$ws.Cells.Item(247, 15).Value = 9000
$ws.Cells.Item(247, 16).Value = 9000
$ws.Cells.Item(247, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(247, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(247, 19).Value = 1286
$ws.Cells.Item(247, 20).Value = 7

# Row 248
$ws.Cells.Item(248, 1).Value = 11
$ws.Cells.Item(248, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(248, 3).Value = "Bíobío"
$ws.Cells.Item(248, 4).Value = 44530
$ws.Cells.Item(248, 5).Value = 8
$ws.Cells.Item(248, 6).Value = "Fruta"
$ws.Cells.Item(248, 7).Value = 100101
$ws.Cells.Item(248, 8).Value = "Berries"
$ws.Cells.Item(248, 9).Value = 100112025
$ws.Cells.Item(248, 10).Value = "Frutilla"
$ws.Cells.Item(248, 11).Value = "Sin especificar"
$ws.Cells.Item(248, 12).Value = "Primera"
$ws.Cells.Item(248, 13).Value = 100
$ws.Cells.Item(248, 14).Value = 7000
$ws.Cells.Item(248, 15).Value = 7000
$ws.Cells.Item(248, 16).Value = 7000
$ws.Cells.Item(248, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(248, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(248, 19).Value = 1000
$ws.Cells.Item(248, 20).Value = 7

# Row 249
$ws.Cells.Item(249, 1).Value = 11
$ws.Cells.Item(249, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(249, 3).Value = "Bíobío"
$ws.Cells.Item(249, 4).Value = 44530
$ws.Cells.Item(249, 5).Value = 8
$ws.Cells.Item(249, 6).Value = "Fruta"
$ws.Cells.Item(249, 7).Value = 100101
$ws.Cells.Item(249, 8).Value = "Berries"
$ws.Cells.Item(249, 9).Value = 100112025
$ws.Cells.Item(249, 10).Value = "Frutilla"
$ws.Cells.Item(249, 11).Value = "Sin especificar"
$ws.Cells.Item(249, 12).Value = "Segunda"
$ws.Cells.Item(249, 13).Value = 100
$ws.Cells.Item(249, 14).Value = 6000
$ws.Cells.Item(249, 15).Value = 6000
$ws.Cells.Item(249, 16).Value = 6000
$ws.Cells.Item(249, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(249, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(249, 19).Value = 857
$ws.Cells.Item(249, 20).Value = 7

# Row 250
$ws.Cells.Item(250, 1).Value = 11
$ws.Cells.Item(250, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(250, 3).Value = "Bíobío"
$ws.Cells.Item(250, 4).Value = 44483
$ws.Cells.Item(250, 5).Value = 8
$ws.Cells.Item(250, 6).Value = "Fruta"
$ws.Cells.Item(250, 7).Value = 100101
$ws.Cells.Item(250, 8).Value = "Berries"
$ws.Cells.Item(250, 9).Value = 100112025
$ws.Cells.Item(250, 10).Value = "Frutilla"
$ws.Cells.Item(250, 11).Value = "Sin especificar"
$ws.Cells.Item(250, 12).Value = "Especial"
$ws.Cells.Item(250, 13).Value = 100
$ws.Cells.Item(250, 14).Value = 10000
$ws.Cells.Item(250, 15).Value = 10000
$ws.Cells.Item(250, 16).Value = 10000
$ws.Cells.Item(250, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(250, 18).Value = "Provincia de Cardenal Caro"
$ws.Cells.Item(250, 19).Value = 1429
$ws.Cells.Item(250, 20).Value = 7

# Row 251
$ws.Cells.Item(251, 1).Value = 11
$ws.Cells.Item(251, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(251, 3).Value = "Bíobío"
$ws.Cells.Item(251, 4).Value = 44483
$ws.Cells.Item(251, 5).Value = 8
$ws.Cells.Item(251, 6).Value = "Fruta"
$ws.Cells.Item(251, 7).Value = 100101
$ws.Cells.Item(251, 8).Value = "Berries"
$ws.Cells.Item(251, 9).Value = 100112025
$ws.Cells.Item(251, 10).Value = "Frutilla"
$ws.Cells.Item(251, 11).Value = "Sin especificar"
$ws.Cells.Item(251, 12).Value = "Primera"
$ws.Cells.Item(251, 13).Value = 100
$ws.Cells.Item(251, 14).Value = 9000
$ws.Cells.Item(251, 15).Value = 9000
$ws.Cells.Item(251, 16).Value = 9000
$ws.Cells.Item(251, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(251, 18).Value = "Provincia de Cardenal Caro"
$ws.Cells.Item(251, 19).Value = 1286
$ws.Cells.Item(251, 20).Value = 7

# Row 252
$ws.Cells.Item(252, 1).Value = 11
$ws.Cells.Item(252, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(252, 3).Value = "Bíobío"
$ws.Cells.Item(252, 4).Value = 44489
$ws.Cells.Item(252, 5).Value = 8
$ws.Cells.Item(252, 6).Value = "Fruta"
$ws.Cells.Item(252, 7).Value = 100101
$ws.Cells.Item(252, 8).Value = "Berries"
$ws.Cells.Item(252, 9).Value = 100112025
$ws.Cells.Item(252, 10).Value = "Frutilla"
$ws.Cells.Item(252, 11).Value = "Sin especificar"
$ws.Cells.Item(252, 12).Value = "Primera"
$ws.Cells.Item(252, 13).Value = 430
$ws.Cells.Item(252, 14).Value = 7500
$ws.Cells.Item(252, 15).Value = 8000
$ws.Cells.Item(252, 16).Value = 7767
$ws.Cells.Item(252, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(252, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(252, 19).Value = 1110
$ws.Cells.Item(252, 20).Value = 7

# Row 253
$ws.Cells.Item(253, 1).Value = 11
$ws.Cells.Item(253, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(253, 3).Value = "Bíobío"
$ws.Cells.Item(253, 4).Value = 44489
$ws.Cells.Item(253, 5).Value = 8
$ws.Cells.Item(253, 6).Value = "Fruta"
$ws.Cells.Item(253, 7).Value = 100101
$ws.Cells.Item(253, 8).Value = "Berries"
$ws.Cells.Item(253, 9).Value = 100112025
$ws.Cells.Item(253, 10).Value = "Frutilla"
$ws.Cells.Item(253, 11).Value = "Sin especificar"
$ws.Cells.Item(253, 12).Value = "Segunda"
$ws.Cells.Item(253, 13).Value = 450
$ws.Cells.Item(253, 14).Value = 6000
$ws.Cells.Item(253, 15).Value = 6500
$ws.Cells.Item(253, 16).Value = 6222
$ws.Cells.Item(253, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(253, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(253, 19).Value = 889
$ws.Cells.Item(253, 20).Value = 7

# Row 254
$ws.Cells.Item(254, 1).Value = 11
$ws.Cells.Item(254, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(254, 3).Value = "Bíobío"
$ws.Cells.Item(254, 4).Value = 44232
$ws.Cells.Item(254, 5).Value = 8
$ws.Cells.Item(254, 6).Value = "Fruta"
$ws.Cells.Item(254, 7).Value = 100101
$ws.Cells.Item(254, 8).Value = "Berries"
$ws.Cells.Item(254, 9).Value = 100112025
$ws.Cells.Item(254, 10).Value = "Frutilla"
$ws.Cells.Item(254, 11).Value = "Sin especificar"
$ws.Cells.Item(254, 12).Value = "Especial"
$ws.Cells.Item(254, 13).Value = 100
$ws.Cells.Item(254, 14).Value = 11000
$ws.Cells.Item(254, 15).Value = 11000
$ws.Cells.Item(254, 16).Value = 11000
$ws.Cells.Item(254, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(254, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(254, 19).Value = 1571
$ws.Cells.Item(254, 20).Value = 7

# Row 255
$ws.Cells.Item(255, 1).Value = 11
$ws.Cells.Item(255, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(255, 3).Value = "Bíobío"
$ws.Cells.Item(255, 4).Value = 44232
$ws.Cells.Item(255, 5).Value = 8
$ws.Cells.Item(255, 6).Value = "Fruta"
$ws.Cells.Item(255, 7).Value = 100101
$ws.Cells.Item(255, 8).Value = "Berries"
$ws.Cells.Item(255, 9).Value = 100112025
$ws.Cells.Item(255, 10).Value = "Frutilla"
$ws.Cells.Item(255, 11).Value = "Sin especificar"
$ws.Cells.Item(255, 12).Value = "Primera"
$ws.Cells.Item(255, 13).Value = 100
$ws.Cells.Item(255, 14).Value = 9000
$ws.Cells.Item(255, 15).Value = 9000
$ws.Cells.Item(255, 16).Value = 9000
$ws.Cells.Item(255, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(255, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(255, 19).Value = 1286
$ws.Cells.Item(255, 20).Value = 7

# Row 256
$ws.Cells.Item(256, 1).Value = 11
$ws.Cells.Item(256, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(256, 3).Value = "Bíobío"
$ws.Cells.Item(256, 4).Value = 44232
$ws.Cells.Item(256, 5).Value = 8
$ws.Cells.Item(256, 6).Value = "Fruta"
$ws.Cells.Item(256, 7).Value = 100101
$ws.Cells.Item(256, 8).Value = "Berries"
$ws.Cells.Item(256, 9).Value = 100112025
$ws.Cells.Item(256, 10).Value = "Frutilla"
$ws.Cells.Item(256, 11).Value = "Sin especificar"
$ws.Cells.Item(256, 12).Value = "Segunda"
$ws.Cells.Item(256, 13).Value = 100
$ws.Cells.Item(256, 14).Value = 7000
$ws.Cells.Item(256, 15).Value = 7000
$ws.Cells.Item(256, 16).Value = 7000
$ws.Cells.Item(256, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(256, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(256, 19).Value = 1000
$ws.Cells.Item(256, 20).Value = 7

# Row 257
$ws.Cells.Item(257, 1).Value = 11
$ws.Cells.Item(257, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(257, 3).Value = "Bíobío"
$ws.Cells.Item(257, 4).Value = 44279
$ws.Cells.Item(257, 5).Value = 8
$ws.Cells.Item(257, 6).Value = "Fruta"
$ws.Cells.Item(257, 7).Value = 100101
$ws.Cells.Item(257, 8).Value = "Berries"
$ws.Cells.Item(257, 9).Value = 100112025
$ws.Cells.Item(257, 10).Value = "Frutilla"
$ws.Cells.Item(257, 11).Value = "Sin especificar"
$ws.Cells.Item(257, 12).Value = "Especial"
$ws.Cells.Item(257, 13).Value = 100
$ws.Cells.Item(257, 14).Value = 11000
$ws.Cells.Item(257, 15).Value = 11000
$ws.Cells.Item(257, 16).Value = 11000
$ws.Cells.Item(257, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(257, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(257, 19).Value = 1571
$ws.Cells.Item(257, 20).Value = 7

# Row 258
$ws.Cells.Item(258, 1).Value = 11
$ws.Cells.Item(258, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(258, 3).Value = "Bíobío"
$ws.Cells.Item(258, 4).Value = 44279
$ws.Cells.Item(258, 5).Value = 8
$ws.Cells.Item(258, 6).Value = "Fruta"
$ws.Cells.Item(258, 7).Value = 100101
$ws.Cells.Item(258, 8).Value = "Berries"
$ws.Cells.Item(258, 9).Value = 100112025
$ws.Cells.Item(258, 10).Value = "Frutilla"
$ws.Cells.Item(258, 11).Value = "Sin especificar"
$ws.Cells.Item(258, 12).Value = "Primera"
$ws.Cells.Item(258, 13).Value = 100
$ws.Cells.Item(258, 14).Value = 9000
$ws.Cells.Item(258, 15).Value = 9000
$ws.Cells.Item(258, 16).Value = 9000
$ws.Cells.Item(258, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(258, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(258, 19).Value = 1286
$ws.Cells.Item(258, 20).Value = 7

# Row 259
$ws.Cells.Item(259, 1).Value = 11
$ws.Cells.Item(259, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(259, 3).Value = "Bíobío"
$ws.Cells.Item(259, 4).Value = 44279
$ws.Cells.Item(259, 5).Value = 8
$ws.Cells.Item(259, 6).Value = "Fruta"
$ws.Cells.Item(259, 7).Value = 100101
$ws.Cells.Item(259, 8).Value = "Berries"
$ws.Cells.Item(259, 9).Value = 100112025
$ws.Cells.Item(259, 10).Value = "Frutilla"
$ws.Cells.Item(259, 11).Value = "Sin especificar"
$ws.Cells.Item(259, 12).Value = "Segunda"
$ws.Cells.Item(259, 13).Value = 50
$ws.Cells.Item(259, 14).Value = 7000
$ws.Cells.Item(259, 15).Value = 7000
$ws.Cells.Item(259, 16).Value = 7000
$ws.Cells.Item(259, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(259, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(259, 19).Value = 1000
$ws.Cells.Item(259, 20).Value = 7

# Row 260
$ws.Cells.Item(260, 1).Value = 11
$ws.Cells.Item(260, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(260, 3).Value = "Bíobío"
$ws.Cells.Item(260, 4).Value = 44504
$ws.Cells.Item(260, 5).Value = 8
$ws.Cells.Item(260, 6).Value = "Fruta"
$ws.Cells.Item(260, 7).Value = 100101
$ws.Cells.Item(260, 8).Value = "Berries"
$ws.Cells.Item(260, 9).Value = 100112025
$ws.Cells.Item(260, 10).Value = "Frutilla"
$ws.Cells.Item(260, 11).Value = "Sin especificar"
$ws.Cells.Item(260, 12).Value = "Especial"
$ws.Cells.Item(260, 13).Value = 100
$ws.Cells.Item(260, 14).Value = 9000
$ws.Cells.Item(260, 15).Value = 9000
$ws.Cells.Item(260, 16).Value = 9000
$ws.Cells.Item(260, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(260, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(260, 19).Value = 1286
$ws.Cells.Item(260, 20).Value = 7

# Row 261
$ws.Cells.Item(261, 1).Value = 11
$ws.Cells.Item(261, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(261, 3).Value = "Bíobío"
$ws.Cells.Item(261, 4).Value = 44504
$ws.Cells.Item(261, 5).Value = 8
$ws.Cells.Item(261, 6).Value = "Fruta"
$ws.Cells.Item(261, 7).Value = 100101
$ws.Cells.Item(261, 8).Value = "Berries"
$ws.Cells.Item(261, 9).Value = 100112025
$ws.Cells.Item(261, 10).Value = "Frutilla"
$ws.Cells.Item(261, 11).Value = "Sin especificar"
$ws.Cells.Item(261, 12).Value = "Primera"
$ws.Cells.Item(261, 13).Value = 100
$ws.Cells.Item(261, 14).Value = 7000
$ws.Cells.Item(261, 15).Value = 7000
$ws.Cells.Item(261, 16).Value = 7000
$ws.Cells.Item(261, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(261, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(261, 19).Value = 1000
$ws.Cells.Item(261, 20).Value = 7

# Row 262
$ws.Cells.Item(262, 1).Value = 11
$ws.Cells.Item(262, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(262, 3).Value = "Bíobío"
$ws.Cells.Item(262, 4).Value = 44504
$ws.Cells.Item(262, 5).Value = 8
$ws.Cells.Item(262, 6).Value = "Fruta"
$ws.Cells.Item(262, 7).Value = 100101
$ws.Cells.Item(262, 8).Value = "Berries"
$ws.Cells.Item(262, 9).Value = 100112025
$ws.Cells.Item(262, 10).Value = "Frutilla"
$ws.Cells.Item(262, 11).Value = "Sin especificar"
$ws.Cells.Item(262, 12).Value = "Segunda"
$ws.Cells.Item(262, 13).Value = 50
$ws.Cells.Item(262, 14).Value = 6000
$ws.Cells.Item(262, 15).Value = 6000
$ws.Cells.Item(262, 16).Value = 6000
$ws.Cells.Item(262, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(262, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(262, 19).Value = 857
$ws.Cells.Item(262, 20).Value = 7

# Row 263
$ws.Cells.Item(263, 1).Value = 11
$ws.Cells.Item(263, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(263, 3).Value = "Bíobío"
$ws.Cells.Item(263, 4).Value = 44572
$ws.Cells.Item(263, 5).Value = 8
$ws.Cells.Item(263, 6).Value = "Fruta"
$ws.Cells.Item(263, 7).Value = 100101
$ws.Cells.Item(263, 8).Value = "Berries"
$ws.Cells.Item(263, 9).Value = 100112025
$ws.Cells.Item(263, 10).Value = "Frutilla"
$ws.Cells.Item(263, 11).Value = "Sin especificar"
$ws.Cells.Item(263, 12).Value = "Especial"
$ws.Cells.Item(263, 13).Value = 200
$ws.Cells.Item(263, 14).Value = 8000
$ws.Cells.Item(263, 15).Value = 8000
$ws.Cells.Item(263, 16).Value = 8000
$ws.Cells.Item(263, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(263, 18).Value = "Región del Maule"
$ws.Cells.Item(263, 19).Value = 1143
$ws.Cells.Item(263, 20).Value = 7

# Row 264
$ws.Cells.Item(264, 1).Value = 11
$ws.Cells.Item(264, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(264, 3).Value = "Bíobío"
$ws.Cells.Item(264, 4).Value = 44572
$ws.Cells.Item(264, 5).Value = 8
$ws.Cells.Item(264, 6).Value = "Fruta"
$ws.Cells.Item(264, 7).Value = 100101
$ws.Cells.Item(264, 8).Value = "Berries"
$ws.Cells.Item(264, 9).Value = 100112025
$ws.Cells.Item(264, 10).Value = "Frutilla"
$ws.Cells.Item(264, 11).Value = "Sin especificar"
$ws.Cells.Item(264, 12).Value = "Primera"
$ws.Cells.Item(264, 13).Value = 200
$ws.Cells.Item(264, 14).Value = 7000
$ws.Cells.Item(264, 15).Value = 7000
$ws.Cells.Item(264, 16).Value = 7000
$ws.Cells.Item(264, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(264, 18).Value = "Región del Maule"
$ws.Cells.Item(264, 19).Value = 1000
$ws.Cells.Item(264, 20).Value = 7

# Row 265
$ws.Cells.Item(265, 1).Value = 11
$ws.Cells.Item(265, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(265, 3).Value = "Bíobío"
$ws.Cells.Item(265, 4).Value = 44572
$ws.Cells.Item(265, 5).Value = 8
$ws.Cells.Item(265, 6).Value = "Fruta"
$ws.Cells.Item(265, 7).Value = 100101
$ws.Cells.Item(265, 8).Value = "Berries"
$ws.Cells.Item(265, 9).Value = 100112025
$ws.Cells.Item(265, 10).Value = "Frutilla"
$ws.Cells.Item(265, 11).Value = "Sin especificar"
$ws.Cells.Item(265, 12).Value = "Segunda"
$ws.Cells.Item(265, 13).Value = 100
$ws.Cells.Item(265, 14).Value = 5000
$ws.Cells.Item(265, 15).Value = 5000
$ws.Cells.Item(265, 16).Value = 5000
$ws.Cells.Item(265, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(265, 18).Value = "Región del Maule"
$ws.Cells.Item(265, 19).Value = 714
$ws.Cells.Item(265, 20).Value = 7

# Row 266
$ws.Cells.Item(266, 1).Value = 11
$ws.Cells.Item(266, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(266, 3).Value = "Bíobío"
$ws.Cells.Item(266, 4).Value = 44257
$ws.Cells.Item(266, 5).Value = 8
$ws.Cells.Item(266, 6).Value = "Fruta"
$ws.Cells.Item(266, 7).Value = 100101
$ws.Cells.Item(266, 8).Value = "Berries"
$ws.Cells.Item(266, 9).Value = 100112025
$ws.Cells.Item(266, 10).Value = "Frutilla"
$ws.Cells.Item(266, 11).Value = "Sin especificar"
$ws.Cells.Item(266, 12).Value = "Primera"
$ws.Cells.Item(266, 13).Value = 500
$ws.Cells.Item(266, 14).Value = 8000
$ws.Cells.Item(266, 15).Value = 10000
$ws.Cells.Item(266, 16).Value = 8800
$ws.Cells.Item(266, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(266, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(266, 19).Value = 1257
$ws.Cells.Item(266, 20).Value = 7

# Row 267
$ws.Cells.Item(267, 1).Value = 11
$ws.Cells.Item(267, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(267, 3).Value = "Bíobío"
$ws.Cells.Item(267, 4).Value = 44257
$ws.Cells.Item(267, 5).Value = 8
$ws.Cells.Item(267, 6).Value = "Fruta"
$ws.Cells.Item(267, 7).Value = 100101
$ws.Cells.Item(267, 8).Value = "Berries"
$ws.Cells.Item(267, 9).Value = 100112025
$ws.Cells.Item(267, 10).Value = "Frutilla"
$ws.Cells.Item(267, 11).Value = "Sin especificar"
$ws.Cells.Item(267, 12).Value = "Segunda"
$ws.Cells.Item(267, 13).Value = 100
$ws.Cells.Item(267, 14).Value = 6000
$ws.Cells.Item(267, 15).Value = 6000
$ws.Cells.Item(267, 16).Value = 6000
$ws.Cells.Item(267, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(267, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(267, 19).Value = 857
$ws.Cells.Item(267, 20).Value = 7

# Row 268
$ws.Cells.Item(268, 1).Value = 11
$ws.Cells.Item(268, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(268, 3).Value = "Bíobío"
$ws.Cells.Item(268, 4).Value = 44385
$ws.Cells.Item(268, 5).Value = 8
$ws.Cells.Item(268, 6).Value = "Fruta"
$ws.Cells.Item(268, 7).Value = 100101
$ws.Cells.Item(268, 8).Value = "Berries"
$ws.Cells.Item(268, 9).Value = 100112025
$ws.Cells.Item(268, 10).Value = "Frutilla"
$ws.Cells.Item(268, 11).Value = "Sin especificar"
$ws.Cells.Item(268, 12).Value = "Especial"
$ws.Cells.Item(268, 13).Value = 100
$ws.Cells.Item(268, 14).Value = 20000
$ws.Cells.Item(268, 15).Value = 20000
$ws.Cells.Item(268, 16).Value = 20000
$ws.Cells.Item(268, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(268, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(268, 19).Value = 2857
$ws.Cells.Item(268, 20).Value = 7

# Row 269
$ws.Cells.Item(269, 1).Value = 11
$ws.Cells.Item(269, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(269, 3).Value = "Bíobío"
$ws.Cells.Item(269, 4).Value = 44385
$ws.Cells.Item(269, 5).Value = 8
$ws.Cells.Item(269, 6).Value = "Fruta"
$ws.Cells.Item(269, 7).Value = 100101
$ws.Cells.Item(269, 8).Value = "Berries"
$ws.Cells.Item(269, 9).Value = 100112025
$ws.Cells.Item(269, 10).Value = "Frutilla"
$ws.Cells.Item(269, 11).Value = "Sin especificar"
$ws.Cells.Item(269, 12).Value = "Segunda"
$ws.Cells.Item(269, 13).Value = 50
$ws.Cells.Item(269, 14).Value = 14000
$ws.Cells.Item(269, 15).Value = 14000
$ws.Cells.Item(269, 16).Value = 14000
$ws.Cells.Item(269, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(269, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(269, 19).Value = 2000
$ws.Cells.Item(269, 20).Value = 7

# Row 270
$ws.Cells.Item(270, 1).Value = 11
$ws.Cells.Item(270, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(270, 3).Value = "Bíobío"
$ws.Cells.Item(270, 4).Value = 44236
$ws.Cells.Item(270, 5).Value = 8
$ws.Cells.Item(270, 6).Value = "Fruta"
$ws.Cells.Item(270, 7).Value = 100101
$ws.Cells.Item(270, 8).Value = "Berries"
$ws.Cells.Item(270, 9).Value = 100112025
$ws.Cells.Item(270, 10).Value = "Frutilla"
$ws.Cells.Item(270, 11).Value = "Sin especificar"
$ws.Cells.Item(270, 12).Value = "Especial"
$ws.Cells.Item(270, 13).Value = 100
$ws.Cells.Item(270, 14).Value = 10000
$ws.Cells.Item(270, 15).Value = 10000
$ws.Cells.Item(270, 16).Value = 10000
$ws.Cells.Item(270, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(270, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(270, 19).Value = 1429
$ws.Cells.Item(270, 20).Value = 7

# Row 271
$ws.Cells.Item(271, 1).Value = 11
$ws.Cells.Item(271, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(271, 3).Value = "Bíobío"
$ws.Cells.Item(271, 4).Value = 44236
$ws.Cells.Item(271, 5).Value = 8
$ws.Cells.Item(271, 6).Value = "Fruta"
$ws.Cells.Item(271, 7).Value = 100101
$ws.Cells.Item(271, 8).Value = "Berries"
$ws.Cells.Item(271, 9).Value = 100112025
$ws.Cells.Item(271, 10).Value = "Frutilla"
$ws.Cells.Item(271, 11).Value = "Sin especificar"
$ws.Cells.Item(271, 12).Value = "Primera"
$ws.Cells.Item(271, 13).Value = 100
$ws.Cells.Item(271, 14).Value = 8000
$ws.Cells.Item(271, 15).Value = 8000
$ws.Cells.Item(271, 16).Value = 8000
$ws.Cells.Item(271, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(271, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(271, 19).Value = 1143
$ws.Cells.Item(271, 20).Value = 7

# Row 272
$ws.Cells.Item(272, 1).Value = 11
$ws.Cells.Item(272, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(272, 3).Value = "Bíobío"
$ws.Cells.Item(272, 4).Value = 44236
$ws.Cells.Item(272, 5).Value = 8
$ws.Cells.Item(272, 6).Value = "Fruta"
$ws.Cells.Item(272, 7).Value = 100101
$ws.Cells.Item(272, 8).Value = "Berries"
$ws.Cells.Item(272, 9).Value = 100112025
$ws.Cells.Item(272, 10).Value = "Frutilla"
$ws.Cells.Item(272, 11).Value = "Sin especificar"
$ws.Cells.Item(272, 12).Value = "Segunda"
$ws.Cells.Item(272, 13).Value = 100
$ws.Cells.Item(272, 14).Value = 7000
$ws.Cells.Item(272, 15).Value = 7000
$ws.Cells.Item(272, 16).Value = 7000
$ws.Cells.Item(272, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(272, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(272, 19).Value = 1000
$ws.Cells.Item(272, 20).Value = 7

# Row 273
$ws.Cells.Item(273, 1).Value = 11
$ws.Cells.Item(273, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(273, 3).Value = "Bíobío"
$ws.Cells.Item(273, 4).Value = 44299
$ws.Cells.Item(273, 5).Value = 8
$ws.Cells.Item(273, 6).Value = "Fruta"
$ws.Cells.Item(273, 7).Value = 100101
$ws.Cells.Item(273, 8).Value = "Berries"
$ws.Cells.Item(273, 9).Value = 100112025
$ws.Cells.Item(273, 10).Value = "Frutilla"
$ws.Cells.Item(273, 11).Value = "Sin especificar"
$ws.Cells.Item(273, 12).Value = "Especial"
$ws.Cells.Item(273, 13).Value = 100
$ws.Cells.Item(273, 14).Value = 11000
$ws.Cells.Item(273, 15).Value = 11000
$ws.Cells.Item(273, 16).Value = 11000
$ws.Cells.Item(273, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(273, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(273, 19).Value = 1571
$ws.Cells.Item(273, 20).Value = 7

# Row 274
$ws.Cells.Item(274, 1).Value = 11
$ws.Cells.Item(274, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(274, 3).Value = "Bíobío"
$ws.Cells.Item(274, 4).Value = 44299
$ws.Cells.Item(274, 5).Value = 8
$ws.Cells.Item(274, 6).Value = "Fruta"
$ws.Cells.Item(274, 7).Value = 100101
$ws.Cells.Item(274, 8).Value = "Berries"
$ws.Cells.Item(274, 9).Value = 100112025
$ws.Cells.Item(274, 10).Value = "Frutilla"
$ws.Cells.Item(274, 11).Value = "Sin especificar"
$ws.Cells.Item(274, 12).Value = "Primera"
$ws.Cells.Item(274, 13).Value = 100
$ws.Cells.Item(274, 14).Value = 9000
$ws.Cells.Item(274, 15).Value = 9000
$ws.Cells.Item(274, 16).Value = 9000
$ws.Cells.Item(274, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(274, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(274, 19).Value = 1286
$ws.Cells.Item(274, 20).Value = 7

# Row 275
$ws.Cells.Item(275, 1).Value = 11
$ws.Cells.Item(275, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(275, 3).Value = "Bíobío"
$ws.Cells.Item(275, 4).Value = 44167
$ws.Cells.Item(275, 5).Value = 8
$ws.Cells.Item(275, 6).Value = "Fruta"
$ws.Cells.Item(275, 7).Value = 100101
$ws.Cells.Item(275, 8).Value = "Berries"
$ws.Cells.Item(275, 9).Value = 100112025
$ws.Cells.Item(275, 10).Value = "Frutilla"
$ws.Cells.Item(275, 11).Value = "Sin especificar"
$ws.Cells.Item(275, 12).Value = "Primera"
$ws.Cells.Item(275, 13).Value = 100
$ws.Cells.Item(275, 14).Value = 8000
$ws.Cells.Item(275, 15).Value = 8000
$ws.Cells.Item(275, 16).Value = 8000
$ws.Cells.Item(275, 17).Value = "$/bandeja 3 kilos"
$ws.Cells.Item(275, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(275, 19).Value = 2667
$ws.Cells.Item(275, 20).Value = 3

# Row 276
$ws.Cells.Item(276, 1).Value = 11
$ws.Cells.Item(276, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(276, 3).Value = "Bíobío"
$ws.Cells.Item(276, 4).Value = 44167
$ws.Cells.Item(276, 5).Value = 8
$ws.Cells.Item(276, 6).Value = "Fruta"
$ws.Cells.Item(276, 7).Value = 100101
$ws.Cells.Item(276, 8).Value = "Berries"
$ws.Cells.Item(276, 9).Value = 100112025
$ws.Cells.Item(276, 10).Value = "Frutilla"
$ws.Cells.Item(276, 11).Value = "Sin especificar"
$ws.Cells.Item(276, 12).Value = "Segunda"
$ws.Cells.Item(276, 13).Value = 50
$ws.Cells.Item(276, 14).Value = 7000
$ws.Cells.Item(276, 15).Value = 7000
$ws.Cells.Item(276, 16).Value = 7000
$ws.Cells.Item(276, 17).Value = "$/bandeja 3 kilos"
$ws.Cells.Item(276, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(276, 19).Value = 2333
$ws.Cells.Item(276, 20).Value = 3

# Row 277
$ws.Cells.Item(277, 1).Value = 11
$ws.Cells.Item(277, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(277, 3).Value = "Bíobío"
$ws.Cells.Item(277, 4).Value = 44349
$ws.Cells.Item(277, 5).Value = 8
$ws.Cells.Item(277, 6).Value = "Fruta"
$ws.Cells.Item(277, 7).Value = 100101
$ws.Cells.Item(277, 8).Value = "Berries"
$ws.Cells.Item(277, 9).Value = 100112025
$ws.Cells.Item(277, 10).Value = "Frutilla"
$ws.Cells.Item(277, 11).Value = "Sin especificar"
$ws.Cells.Item(277, 12).Value = "Especial"
$ws.Cells.Item(277, 13).Value = 50
$ws.Cells.Item(277, 14).Value = 14000
$ws.Cells.Item(277, 15).Value = 14000
$ws.Cells.Item(277, 16).Value = 14000
$ws.Cells.Item(277, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(277, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(277, 19).Value = 2000
$ws.Cells.Item(277, 20).Value = 7

# Row 278
$ws.Cells.Item(278, 1).Value = 11
$ws.Cells.Item(278, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(278, 3).Value = "Bíobío"
$ws.Cells.Item(278, 4).Value = 44349
$ws.Cells.Item(278, 5).Value = 8
$ws.Cells.Item(278, 6).Value = "Fruta"
$ws.Cells.Item(278, 7).Value = 100101
$ws.Cells.Item(278, 8).Value = "Berries"
$ws.Cells.Item(278, 9).Value = 100112025
$ws.Cells.Item(278, 10).Value = "Frutilla"
$ws.Cells.Item(278, 11).Value = "Sin especificar"
$ws.Cells.Item(278, 12).Value = "Primera"
$ws.Cells.Item(278, 13).Value = 50
$ws.Cells.Item(278, 14).Value = 11000
$ws.Cells.Item(278, 15).Value = 11000
$ws.Cells.Item(278, 16).Value = 11000
$ws.Cells.Item(278, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(278, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(278, 19).Value = 1571
$ws.Cells.Item(278, 20).Value = 7

# Row 279
$ws.Cells.Item(279, 1).Value = 11
$ws.Cells.Item(279, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(279, 3).Value = "Bíobío"
$ws.Cells.Item(279, 4).Value = 44285
$ws.Cells.Item(279, 5).Value = 8
$ws.Cells.Item(279, 6).Value = "Fruta"
$ws.Cells.Item(279, 7).Value = 100101
$ws.Cells.Item(279, 8).Value = "Berries"
$ws.Cells.Item(279, 9).Value = 100112025
$ws.Cells.Item(279, 10).Value = "Frutilla"
$ws.Cells.Item(279, 11).Value = "Sin especificar"
$ws.Cells.Item(279, 12).Value = "Primera"
$ws.Cells.Item(279, 13).Value = 100
$ws.Cells.Item(279, 14).Value = 10000
$ws.Cells.Item(279, 15).Value = 10000
$ws.Cells.Item(279, 16).Value = 10000
$ws.Cells.Item(279, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(279, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(279, 19).Value = 1429
$ws.Cells.Item(279, 20).Value = 7

# Row 280
$ws.Cells.Item(280, 1).Value = 11
$ws.Cells.Item(280, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(280, 3).Value = "Bíobío"
$ws.Cells.Item(280, 4).Value = 44285
$ws.Cells.Item(280, 5).Value = 8
$ws.Cells.Item(280, 6).Value = "Fruta"
$ws.Cells.Item(280, 7).Value = 100101
$ws.Cells.Item(280, 8).Value = "Berries"
$ws.Cells.Item(280, 9).Value = 100112025
$ws.Cells.Item(280, 10).Value = "Frutilla"
$ws.Cells.Item(280, 11).Value = "Sin especificar"
$ws.Cells.Item(280, 12).Value = "Segunda"
$ws.Cells.Item(280, 13).Value = 100
$ws.Cells.Item(280, 14).Value = 8000
$ws.Cells.Item(280, 15).Value = 8000
$ws.Cells.Item(280, 16).Value = 8000
$ws.Cells.Item(280, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(280, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(280, 19).Value = 1143
$ws.Cells.Item(280, 20).Value = 7

# Row 281
$ws.Cells.Item(281, 1).Value = 11
$ws.Cells.Item(281, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(281, 3).Value = "Bíobío"
$ws.Cells.Item(281, 4).Value = 44498
$ws.Cells.Item(281, 5).Value = 8
$ws.Cells.Item(281, 6).Value = "Fruta"
$ws.Cells.Item(281, 7).Value = 100101
$ws.Cells.Item(281, 8).Value = "Berries"
$ws.Cells.Item(281, 9).Value = 100112025
$ws.Cells.Item(281, 10).Value = "Frutilla"
$ws.Cells.Item(281, 11).Value = "Sin especificar"
$ws.Cells.Item(281, 12).Value = "Especial"
$ws.Cells.Item(281, 13).Value = 100
$ws.Cells.Item(281, 14).Value = 10000
$ws.Cells.Item(281, 15).Value = 10000
$ws.Cells.Item(281, 16).Value = 10000
$ws.Cells.Item(281, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(281, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(281, 19).Value = 1429
$ws.Cells.Item(281, 20).Value = 7

# Row 282
$ws.Cells.Item(282, 1).Value = 11
$ws.Cells.Item(282, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(282, 3).Value = "Bíobío"
$ws.Cells.Item(282, 4).Value = 44498
$ws.Cells.Item(282, 5).Value = 8
$ws.Cells.Item(282, 6).Value = "Fruta"
$ws.Cells.Item(282, 7).Value = 100101
$ws.Cells.Item(282, 8).Value = "Berries"
$ws.Cells.Item(282, 9).Value = 100112025
$ws.Cells.Item(282, 10).Value = "Frutilla"
$ws.Cells.Item(282, 11).Value = "Sin especificar"
$ws.Cells.Item(282, 12).Value = "Primera"
$ws.Cells.Item(282, 13).Value = 100
$ws.Cells.Item(282, 14).Value = 8000
$ws.Cells.Item(282, 15).Value = 8000
$ws.Cells.Item(282, 16).Value = 8000
$ws.Cells.Item(282, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(282, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(282, 19).Value = 1143
$ws.Cells.Item(282, 20).Value = 7

# Row 283
$ws.Cells.Item(283, 1).Value = 11
$ws.Cells.Item(283, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(283, 3).Value = "Bíobío"
$ws.Cells.Item(283, 4).Value = 44498
$ws.Cells.Item(283, 5).Value = 8
$ws.Cells.Item(283, 6).Value = "Fruta"
$ws.Cells.Item(283, 7).Value = 100101
$ws.Cells.Item(283, 8).Value = "Berries"
$ws.Cells.Item(283, 9).Value = 100112025
$ws.Cells.Item(283, 10).Value = "Frutilla"
$ws.Cells.Item(283, 11).Value = "Sin especificar"
$ws.Cells.Item(283, 12).Value = "Segunda"
$ws.Cells.Item(283, 13).Value = 50
$ws.Cells.Item(283, 14).Value = 6000
$ws.Cells.Item(283, 15).Value = 6000
$ws.Cells.Item(283, 16).Value = 6000
$ws.Cells.Item(283, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(283, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(283, 19).Value = 857
$ws.Cells.Item(283, 20).Value = 7

# Row 284
$ws.Cells.Item(284, 1).Value = 11
$ws.Cells.Item(284, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(284, 3).Value = "Bíobío"
$ws.Cells.Item(284, 4).Value = 44595
$ws.Cells.Item(284, 5).Value = 8
$ws.Cells.Item(284, 6).Value = "Fruta"
$ws.Cells.Item(284, 7).Value = 100101
$ws.Cells.Item(284, 8).Value = "Berries"
$ws.Cells.Item(284, 9).Value = 100112025
$ws.Cells.Item(284, 10).Value = "Frutilla"
$ws.Cells.Item(284, 11).Value = "Sin especificar"
$ws.Cells.Item(284, 12).Value = "Primera"
$ws.Cells.Item(284, 13).Value = 250
$ws.Cells.Item(284, 14).Value = 6500
$ws.Cells.Item(284, 15).Value = 7000
$ws.Cells.Item(284, 16).Value = 6740
$ws.Cells.Item(284, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(284, 18).Value = "Región del Maule"
$ws.Cells.Item(284, 19).Value = 963
$ws.Cells.Item(284, 20).Value = 7

# Row 285
$ws.Cells.Item(285, 1).Value = 11
$ws.Cells.Item(285, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(285, 3).Value = "Bíobío"
$ws.Cells.Item(285, 4).Value = 44335
$ws.Cells.Item(285, 5).Value = 8
$ws.Cells.Item(285, 6).Value = "Fruta"
$ws.Cells.Item(285, 7).Value = 100101
$ws.Cells.Item(285, 8).Value = "Berries"
$ws.Cells.Item(285, 9).Value = 100112025
$ws.Cells.Item(285, 10).Value = "Frutilla"
$ws.Cells.Item(285, 11).Value = "Sin especificar"
$ws.Cells.Item(285, 12).Value = "Especial"
$ws.Cells.Item(285, 13).Value = 100
$ws.Cells.Item(285, 14).Value = 14000
$ws.Cells.Item(285, 15).Value = 14000
$ws.Cells.Item(285, 16).Value = 14000
$ws.Cells.Item(285, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(285, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(285, 19).Value = 2000
$ws.Cells.Item(285, 20).Value = 7

# Row 286
$ws.Cells.Item(286, 1).Value = 11
$ws.Cells.Item(286, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(286, 3).Value = "Bíobío"
$ws.Cells.Item(286, 4).Value = 44335
$ws.Cells.Item(286, 5).Value = 8
$ws.Cells.Item(286, 6).Value = "Fruta"
$ws.Cells.Item(286, 7).Value = 100101
$ws.Cells.Item(286, 8).Value = "Berries"
$ws.Cells.Item(286, 9).Value = 100112025
$ws.Cells.Item(286, 10).Value = "Frutilla"
$ws.Cells.Item(286, 11).Value = "Sin especificar"
$ws.Cells.Item(286, 12).Value = "Primera"
$ws.Cells.Item(286, 13).Value = 100
$ws.Cells.Item(286, 14).Value = 12000
$ws.Cells.Item(286, 15).Value = 12000
$ws.Cells.Item(286, 16).Value = 12000
$ws.Cells.Item(286, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(286, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(286, 19).Value = 1714
$ws.Cells.Item(286, 20).Value = 7

# Row 287
$ws.Cells.Item(287, 1).Value = 11
$ws.Cells.Item(287, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(287, 3).Value = "Bíobío"
$ws.Cells.Item(287, 4).Value = 44552
$ws.Cells.Item(287, 5).Value = 8
$ws.Cells.Item(287, 6).Value = "Fruta"
$ws.Cells.Item(287, 7).Value = 100101
$ws.Cells.Item(287, 8).Value = "Berries"
$ws.Cells.Item(287, 9).Value = 100112025
$ws.Cells.Item(287, 10).Value = "Frutilla"
$ws.Cells.Item(287, 11).Value = "Sin especificar"
$ws.Cells.Item(287, 12).Value = "Especial"
$ws.Cells.Item(287, 13).Value = 390
$ws.Cells.Item(287, 14).Value = 7500
$ws.Cells.Item(287, 15).Value = 8000
$ws.Cells.Item(287, 16).Value = 7756
$ws.Cells.Item(287, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(287, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(287, 19).Value = 1108
$ws.Cells.Item(287, 20).Value = 7

# Row 288
$ws.Cells.Item(288, 1).Value = 11
$ws.Cells.Item(288, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(288, 3).Value = "Bíobío"
$ws.Cells.Item(288, 4).Value = 44552
$ws.Cells.Item(288, 5).Value = 8
$ws.Cells.Item(288, 6).Value = "Fruta"
$ws.Cells.Item(288, 7).Value = 100101
$ws.Cells.Item(288, 8).Value = "Berries"
$ws.Cells.Item(288, 9).Value = 100112025
$ws.Cells.Item(288, 10).Value = "Frutilla"
$ws.Cells.Item(288, 11).Value = "Sin especificar"
$ws.Cells.Item(288, 12).Value = "Primera"
$ws.Cells.Item(288, 13).Value = 220
$ws.Cells.Item(288, 14).Value = 6000
$ws.Cells.Item(288, 15).Value = 6500
$ws.Cells.Item(288, 16).Value = 6227
$ws.Cells.Item(288, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(288, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(288, 19).Value = 890
$ws.Cells.Item(288, 20).Value = 7

# Row 289
$ws.Cells.Item(289, 1).Value = 11
$ws.Cells.Item(289, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(289, 3).Value = "Bíobío"
$ws.Cells.Item(289, 4).Value = 44544
$ws.Cells.Item(289, 4).NumberFormat = $dateFormat
$ws.Cells.Item(289, 5).Value = 8
$ws.Cells.Item(289, 6).Value = "Fruta"
$ws.Cells.Item(289, 7).Value = 100101
$ws.Cells.Item(289, 8).Value = "Berries"
$ws.Cells.Item(289, 9).Value = 100112025
$ws.Cells.Item(289, 10).Value = "Frutilla"
$ws.Cells.Item(289, 11).Value = "Sin especificar"
$ws.Cells.Item(289, 12).Value = "Especial"
$ws.Cells.Item(289, 13).Value = 150
$ws.Cells.Item(289, 14).Value = 8000
$ws.Cells.Item(289, 15).Value = 8000
$ws.Cells.Item(289, 16).Value = 8000
$ws.Cells.Item(289, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(289, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(289, 19).Value = 1143
$ws.Cells.Item(289, 20).Value = 7

# Row 290
$ws.Cells.Item(290, 1).Value = 11
$ws.Cells.Item(290, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(290, 3).Value = "Bíobío"
$ws.Cells.Item(290, 4).Value = 44544
$ws.Cells.Item(290, 4).NumberFormat = $dateFormat
$ws.Cells.Item(290, 5).Value = 8
$ws.Cells.Item(290, 6).Value = "Fruta"
$ws.Cells.Item(290, 7).Value = 100101
$ws.Cells.Item(290, 8).Value = "Berries"
$ws.Cells.Item(290, 9).Value = 100112025
$ws.Cells.Item(290, 10).Value = "Frutilla"
$ws.Cells.Item(290, 11).Value = "Sin especificar"
$ws.Cells.Item(290, 12).Value = "Primera"
$ws.Cells.Item(290, 13).Value = 100
$ws.Cells.Item(290, 14).Value = 7000
$ws.Cells.Item(290, 15).Value = 7000
$ws.Cells.Item(290, 16).Value = 7000
$ws.Cells.Item(290, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(290, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(290, 19).Value = 1000
$ws.Cells.Item(290, 20).Value = 7

# Row 291
$ws.Cells.Item(291, 1).Value = 11
$ws.Cells.Item(291, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(291, 3).Value = "Bíobío"
$ws.Cells.Item(291, 4).Value = 44544
$ws.Cells.Item(291, 4).NumberFormat = $dateFormat
$ws.Cells.Item(291, 5).Value = 8
$ws.Cells.Item(291, 6).Value = "Fruta"
$ws.Cells.Item(291, 7).Value = 100101
$ws.Cells.Item(291, 8).Value = "Berries"
$ws.Cells.Item(291, 9).Value = 100112025
$ws.Cells.Item(291, 10).Value = "Frutilla"
$ws.Cells.Item(291, 11).Value = "Sin especificar"
$ws.Cells.Item(291, 12).Value = "Segunda"
$ws.Cells.Item(291, 13).Value = 50
$ws.Cells.Item(291, 14).Value = 5000
$ws.Cells.Item(291, 15).Value = 5000
$ws.Cells.Item(291, 16).Value = 5000
$ws.Cells.Item(291, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(291, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(291, 19).Value = 714
$ws.Cells.Item(291, 20).Value = 7
